$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 11076
$ws.Range("E2").Value = 469
$ws.Range("F2").Value = 469
$ws.Range("G2").Value = 413
$ws.Range("H2").Value = 334
$ws.Range("I2").Value = 333
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5717
$ws.Range("L2").Value = 3482
$ws.Range("M2").Value = 2236
$ws.Range("N2").Value = 2238
$ws.Range("O2").Value = -2
$ws.Range("P2").Value = 431
$ws.Range("Q2").Value = 424
$ws.Range("R2").Value = -206
$ws.Range("S2").Value = -199
$ws.Range("T2").Value = 284
$ws.Range("U2").Value = 140
$ws.Range("V2").Value = 1832
$ws.Range("W2").Value = 4.23
$ws.Range("X2").Value = 3.01
$ws.Range("Y2").Value = 15.86
$ws.Range("Z2").Value = 5.92
$ws.Range("AA2").Value = 155.74
$ws.Range("AB2").Value = 391.9
$ws.Range("AC2").Value = 3861
$ws.Range("AD2").Value = 38.85
$ws.Range("AE2").Value = 26734
$ws.Range("AF2").Value = 5.61
$ws.Range("AG2").Value = 762
$ws.Range("AH2").Value = 0.51
$ws.Range("AI2").Value = 12.9
$ws.Range("AJ2").Value = 8629009

# Row 3
$ws.Range("D3").Value = 13738
$ws.Range("E3").Value = 569
$ws.Range("F3").Value = 569
$ws.Range("G3").Value = 516
$ws.Range("H3").Value = 376
$ws.Range("I3").Value = 376
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 6129
$ws.Range("L3").Value = 3591
$ws.Range("M3").Value = 2538
$ws.Range("N3").Value = 2540
$ws.Range("O3").Value = -2
$ws.Range("P3").Value = 431
$ws.Range("Q3").Value = 443
$ws.Range("R3").Value = -457
$ws.Range("S3").Value = 4
$ws.Range("T3").Value = 413
$ws.Range("U3").Value = 30
$ws.Range("V3").Value = 1880
$ws.Range("W3").Value = 4.14
$ws.Range("X3").Value = 2.74
$ws.Range("Y3").Value = 15.73
$ws.Range("Z3").Value = 6.35
$ws.Range("AA3").Value = 141.53
$ws.Range("AB3").Value = 461.75
$ws.Range("AC3").Value = 4354
$ws.Range("AD3").Value = 63.5
$ws.Range("AE3").Value = 30335
$ws.Range("AF3").Value = 9.11
$ws.Range("AG3").Value = 869
$ws.Range("AH3").Value = 0.31
$ws.Range("AI3").Value = 13.04
$ws.Range("AJ3").Value = 8629009

# Row 4
$ws.Range("D4").Value = 18703
$ws.Range("E4").Value = 655
$ws.Range("F4").Value = 655
$ws.Range("G4").Value = 618
$ws.Range("H4").Value = 495
$ws.Range("I4").Value = 495
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6746
$ws.Range("L4").Value = 3823
$ws.Range("M4").Value = 2923
$ws.Range("N4").Value = 2925
$ws.Range("O4").Value = -2
$ws.Range("P4").Value = 431
$ws.Range("Q4").Value = 796
$ws.Range("R4").Value = -506
$ws.Range("S4").Value = -346
$ws.Range("T4").Value = 280
$ws.Range("U4").Value = 517
$ws.Range("V4").Value = 1627
$ws.Range("W4").Value = 3.5
$ws.Range("X4").Value = 2.65
$ws.Range("Y4").Value = 18.12
$ws.Range("Z4").Value = 7.69
$ws.Range("AA4").Value = 130.78
$ws.Range("AB4").Value = 561.37
$ws.Range("AC4").Value = 5737
$ws.Range("AD4").Value = 29.54
$ws.Range("AE4").Value = 35095
$ws.Range("AF4").Value = 4.83
$ws.Range("AG4").Value = 956
$ws.Range("AH4").Value = 0.56
$ws.Range("AI4").Value = 10.83
$ws.Range("AJ4").Value = 8629009

# Row 5
$ws.Range("D5").Value = 20655
$ws.Range("E5").Value = 547
$ws.Range("F5").Value = 547
$ws.Range("G5").Value = 505
$ws.Range("H5").Value = 382
$ws.Range("I5").Value = 382
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 7344
$ws.Range("L5").Value = 4544
$ws.Range("M5").Value = 2800
$ws.Range("N5").Value = 2802
$ws.Range("O5").Value = -2
$ws.Range("P5").Value = 431
$ws.Range("Q5").Value = 631
$ws.Range("R5").Value = -509
$ws.Range("S5").Value = 9
$ws.Range("T5").Value = 580
$ws.Range("U5").Value = 51
$ws.Range("V5").Value = 2137
$ws.Range("W5").Value = 2.65
$ws.Range("X5").Value = 1.85
$ws.Range("Y5").Value = 13.33
$ws.Range("Z5").Value = 5.42
$ws.Range("AA5").Value = 162.27
$ws.Range("AB5").Value = 636.52
$ws.Range("AC5").Value = 4424
$ws.Range("AD5").Value = 35.38
$ws.Range("AE5").Value = 34567
$ws.Range("AF5").Value = 4.53
$ws.Range("AG5").Value = 956
$ws.Range("AH5").Value = 0.61
$ws.Range("AI5").Value = 13.45
$ws.Range("AJ5").Value = 8629009

# Row 6
$ws.Range("D6").Value = 22009
$ws.Range("E6").Value = 599
$ws.Range("F6").Value = 599
$ws.Range("G6").Value = 531
$ws.Range("H6").Value = 422
$ws.Range("I6").Value = 421
$ws.Range("K6").Value = 7523
$ws.Range("L6").Value = 4383
$ws.Range("M6").Value = 3140
$ws.Range("N6").Value = 3141
$ws.Range("P6").Value = 431
$ws.Range("Q6").Value = 683
$ws.Range("R6").Value = -437
$ws.Range("S6").Value = -364
$ws.Range("T6").Value = 414
$ws.Range("U6").Value = 269
$ws.Range("V6").Value = 1824
$ws.Range("W6").Value = 2.72
$ws.Range("X6").Value = 1.92
$ws.Range("Y6").Value = 14.18
$ws.Range("Z6").Value = 5.67
$ws.Range("AA6").Value = 139.61
$ws.Range("AB6").Value = 715.1
$ws.Range("AC6").Value = 4883
$ws.Range("AD6").Value = 26.11
$ws.Range("AE6").Value = 38749
$ws.Range("AF6").Value = 3.29
$ws.Range("AG6").Value = 1004
$ws.Range("AH6").Value = 0.79
$ws.Range("AI6").Value = 12.73
$ws.Range("AJ6").Value = 8629009

# Row 7
$ws.Range("D7").Value = 24422
$ws.Range("E7").Value = 522
$ws.Range("G7").Value = 422
$ws.Range("H7").Value = 333
$ws.Range("I7").Value = 344
$ws.Range("K7").Value = 8785
$ws.Range("L7").Value = 5369
$ws.Range("M7").Value = 3416
$ws.Range("N7").Value = 3415
$ws.Range("P7").Value = 430
$ws.Range("Q7").Value = 905
$ws.Range("R7").Value = -1121
$ws.Range("S7").Value = -83
$ws.Range("T7").Value = 656
$ws.Range("U7").Value = 161
$ws.Range("W7").Value = 2.14
$ws.Range("X7").Value = 1.37
$ws.Range("Y7").Value = 10.48
$ws.Range("Z7").Value = 4.09
$ws.Range("AA7").Value = 157.17
$ws.Range("AC7").Value = 3983
$ws.Range("AD7").Value = 19.13
$ws.Range("AE7").Value = 42133
$ws.Range("AF7").Value = 1.81
$ws.Range("AG7").Value = 1052
$ws.Range("AH7").Value = 1.38
$ws.Range("AI7").Value = 26.4

# Row 8
$ws.Range("D8").Value = 26221
$ws.Range("E8").Value = 703
$ws.Range("G8").Value = 629
$ws.Range("H8").Value = 505
$ws.Range("I8").Value = 507
$ws.Range("K8").Value = 9210
$ws.Range("L8").Value = 5409
$ws.Range("M8").Value = 3801
$ws.Range("N8").Value = 3808
$ws.Range("P8").Value = 430
$ws.Range("Q8").Value = 922
$ws.Range("R8").Value = -638
$ws.Range("S8").Value = -96
$ws.Range("T8").Value = 624
$ws.Range("U8").Value = 299
$ws.Range("W8").Value = 2.68
$ws.Range("X8").Value = 1.93
$ws.Range("Y8").Value = 14.04
$ws.Range("Z8").Value = 5.62
$ws.Range("AA8").Value = 142.32
$ws.Range("AC8").Value = 5874
$ws.Range("AD8").Value = 12.97
$ws.Range("AE8").Value = 46976
$ws.Range("AF8").Value = 1.62
$ws.Range("AG8").Value = 1122
$ws.Range("AH8").Value = 1.47
$ws.Range("AI8").Value = 19.1

# Row 9
$ws.Range("D9").Value = 27834
$ws.Range("E9").Value = 796
$ws.Range("G9").Value = 715
$ws.Range("H9").Value = 569
$ws.Range("I9").Value = 581
$ws.Range("K9").Value = 9690
$ws.Range("L9").Value = 5416
$ws.Range("M9").Value = 4274
$ws.Range("N9").Value = 4302
$ws.Range("P9").Value = 430
$ws.Range("Q9").Value = 1002
$ws.Range("R9").Value = -671
$ws.Range("S9").Value = -105
$ws.Range("T9").Value = 660
$ws.Range("U9").Value = 337
$ws.Range("W9").Value = 2.86
$ws.Range("X9").Value = 2.04
$ws.Range("Y9").Value = 14.32
$ws.Range("Z9").Value = 6.02
$ws.Range("AA9").Value = 126.73
$ws.Range("AC9").Value = 6731
$ws.Range("AD9").Value = 11.32
$ws.Range("AE9").Value = 53083
$ws.Range("AF9").Value = 1.44
$ws.Range("AG9").Value = 1183
$ws.Range("AH9").Value = 1.55
$ws.Range("AI9").Value = 17.57

Write-Host "Applied all IFRS list corrections"
